$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "-"
$ws.Range("E3").Value = "[-, 'MCT-3A-CAM', -, -]"
$ws.Range("E4").Value = "[-, 'MCT-3A-CAM', -, -]"
$ws.Range("B6").Value = "-"
$ws.Range("E6").Value = "[-, 'MCT-3A-CAM', -, -]"
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "[-, 'MCT-3A-CAM', -, -]"
